# Apply "Add data for 2022-04-16" update:
#  - Rename sheet "Through 2022-04-07" -> "Through 2022-04-08"
#  - Update header text "2022 (through 04-07)" -> "2022 (through 04-08)"
#  - Update I5 (May row, current-year column) 24 -> 27
#  - Update I14 (Total row, current-year column) 458 -> 461

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "Through 2022-04-08"

# Update the header cell text that shows the "through" date
$ws.Range("I1").Value = "2022 (through 04-08)"

# Update the updated data point(s)
$ws.Range("I5").Value = 27
$ws.Range("I14").Value = 461
